$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-533 all currently hold the serial date
# value 45179 (2023-09-10) and must be updated to 45180 (2023-09-11).
$ws.Range("C2:C533").Value = 45180
